# Actualización desde MV -datos-
# Adds the new quarterly data row (01-04-2021) to the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 76

# Column A holds a date-like label ("01-04-2021") that must be stored as
# literal text (shared string), not auto-converted to a date serial number.
# Writing the string straight to .Value makes Excel's smart-entry parser
# recognise it as a date (and stamp a date NumberFormat on the cell), so we
# build the text with a formula in a scratch cell (formula results are never
# re-interpreted) and paste just the resulting value back - that keeps the
# cell's style untouched (no NumberFormat residue) while still landing the
# text as a normal shared string.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.Formula = "=""01-04-2021"""
$scratch.Copy()
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.PasteSpecial(-4163)
$scratch.Clear()

$values = @(71795, 439, 2975, 886, 2089, 38730, 38730, 0, 29380, 29055, 325, 188, 188, 83, -76, 71871, 40315, 39457, 858, 31889, 31292, 597, 857, 0, 857, -1225, -1225, 35)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item($row, $col).Value = $values[$i]
}
